# Updates the cryptos list (coin price / 1h volume%) table on Sheet1,
# rows 2-51, columns B (Coin), C (Link), D (Price), E (Volume(1h)).
#
# D-column price strings are written with a leading apostrophe and then
# have their style reset to "Normal" - this forces Excel to keep them as
# literal text (matching the workbook's original inlineStr cells) instead
# of auto-coercing them into numbers, which would corrupt values such as
# "43.949.61" (two dots - not a real number) or drop a trailing zero in
# "2.50" -> 2.5. E-column percentages keep their double-space padding and
# are never numeric-looking enough for Excel to coerce, so they are set
# directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'43.949.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.04%  "
# Row 3
$ws.Range("D3").Value = "'2.358.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.13%  "
# Row 4
$ws.Range("E4").Value = "  +0.07%  "
# Row 5
$ws.Range("D5").Value = "'240.03"
$ws.Range("D5").Style = "Normal"
# Row 7
$ws.Range("D7").Value = "'73.37"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.54%  "
# Row 8
$ws.Range("E8").Value = "  +0.08%  "
# Row 9
$ws.Range("D9").Value = "'0.608"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.51%  "
# Row 10
$ws.Range("E10").Value = "  +1.47%  "
# Row 11
$ws.Range("D11").Value = "'60.63"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.00%  "
# Row 12
$ws.Range("D12").Value = "'33.97"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.74%  "
# Row 13
$ws.Range("D13").Value = "'0.108"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.39%  "
# Row 14
$ws.Range("D14").Value = "'7.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.04%  "
# Row 15
$ws.Range("D15").Value = "'16.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.79%  "
# Row 16
$ws.Range("D16").Value = "'0.906"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.20%  "
# Row 17
$ws.Range("D17").Value = "'2.357.38"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.22%  "
# Row 18
$ws.Range("D18").Value = "'43.892.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.04%  "
# Row 19
$ws.Range("E19").Value = "  +0.81%  "
# Row 20
$ws.Range("D20").Value = "'77.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.81%  "
# Row 21
$ws.Range("D21").Value = "'6.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.51%  "
# Row 22
$ws.Range("D22").Value = "'252.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.43%  "
# Row 23
$ws.Range("E23").Value = "  +0.01%  "
# Row 24
$ws.Range("E24").Value = "  +2.78%  "
# Row 25
$ws.Range("E25").Value = "  -5.89%  "
# Row 26
$ws.Range("D26").Value = "'2.50"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.17%  "
# Row 27
$ws.Range("E27").Value = "  -2.41%  "
# Row 28
$ws.Range("E28").Value = "  +1.65%  "
# Row 29
$ws.Range("D29").Value = "'176.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.74%  "
# Row 30
$ws.Range("D30").Value = "'22.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.14%  "
# Row 31
$ws.Range("D31").Value = "'0.128"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.50%  "
# Row 32
$ws.Range("E32").Value = "  -2.19%  "
# Row 33
$ws.Range("E33").Value = "  -1.49%  "
# Row 34
$ws.Range("E34").Value = "  -3.43%  "
# Row 35
$ws.Range("E35").Value = "  -1.83%  "
# Row 36
$ws.Range("D36").Value = "'3.79"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.10%  "
# Row 37
$ws.Range("D37").Value = "'6.62"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.89%  "
# Row 38
$ws.Range("E38").Value = "  +1.67%  "
# Row 39
$ws.Range("E39").Value = "  -1.30%  "
# Row 40
$ws.Range("D40").Value = "'5.47"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +15.63%  "
# Row 41
$ws.Range("D41").Value = "'64.82"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +11.57%  "
# Row 42
$ws.Range("D42").Value = "'19.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.22%  "
# Row 43
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").Value = "'0.106"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.52%  "
# Row 44
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'9.04"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.26%  "
# Row 45
$ws.Range("E45").Value = "  -2.09%  "
# Row 46
$ws.Range("E46").Value = "  -0.08%  "
# Row 47
$ws.Range("E47").Value = "  -1.03%  "
# Row 48
$ws.Range("E48").Value = "  -2.12%  "
# Row 49
$ws.Range("E49").Value = "  -2.00%  "
# Row 50
$ws.Range("D50").Value = "'98.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.16%  "
# Row 51
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").Value = "'2.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.83%  "
